# Add a new "Program" worksheet (programName / programDesc lookup data)
# to the workbook, as the last sheet, and make it the active tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (BatchSheet) so it
# lands at the end of the tab strip, matching the authored workbook.
$lastSheet   = $wb.Worksheets.Item($wb.Worksheets.Count)
$programSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$programSheet.Name = "Program"

# Header row.
$programSheet.Range("A1").Value = "programName"
$programSheet.Range("B1").Value = "programDesc"

# Data row - write B2 before A2 so the new shared strings land in the
# same order ("Team7LMS" then "Team7LMS18") as the authored file.
$programSheet.Range("B2").Value = "Team7LMS"
$programSheet.Range("A2").Value = "Team7LMS18"

# Cosmetic column widths to match the authored sheet.
$programSheet.Columns.Item(1).ColumnWidth = 17.16666666666667
$programSheet.Columns.Item(2).ColumnWidth = 25.5

# Leave the cursor on A2, as in the authored workbook.
$programSheet.Range("A2").Select() | Out-Null
